$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.441.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.736.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.51%  "

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4619"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +8.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3522"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.08%  "

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07353"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.71%  "

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.55"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.077"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.06%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.38"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.44%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.045"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.740.37"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.94"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001051"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06355"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.61"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.724"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.507.27"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.06"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.098"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.38"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.80"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.934.76"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.02%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "124.42"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.034"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -4.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.044"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09131"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.663"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.389"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.04%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.57"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05979"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2062"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.891"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.62%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.176"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.373"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.704"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.01"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.693"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5791"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.81"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.915"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06825"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.114"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.11"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.58%  "
